# Auto-generated: update crypto price/volume table cells per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'26.520.88"
$ws.Range("E2").Value2 = "  +1.01%  "
$ws.Range("D3").Value2 = "'1.727.80"
$ws.Range("E3").Value2 = "  +0.53%  "
$ws.Range("D4").Value2 = "'0.9992"
$ws.Range("E4").Value2 = "  -0.04%  "
$ws.Range("D5").Value2 = "'244.41"
$ws.Range("E5").Value2 = "  +2.03%  "
$ws.Range("D6").Value2 = "'0.9996"
$ws.Range("E6").Value2 = "  -0.07%  "
$ws.Range("D7").Value2 = "'0.4810"
$ws.Range("E7").Value2 = "  +1.94%  "
$ws.Range("E8").Value2 = "  +1.89%  "
$ws.Range("D9").Value2 = "'0.06184"
$ws.Range("E9").Value2 = "  -0.18%  "
$ws.Range("D10").Value2 = "'1.734.38"
$ws.Range("E10").Value2 = "  +0.92%  "
$ws.Range("D11").Value2 = "'0.07189"
$ws.Range("E11").Value2 = "  +1.65%  "
$ws.Range("E12").Value2 = "  +1.40%  "
$ws.Range("D13").Value2 = "'0.6118"
$ws.Range("E13").Value2 = "  +2.69%  "
$ws.Range("E14").Value2 = "  +2.31%  "
$ws.Range("D15").Value2 = "'77.17"
$ws.Range("E15").Value2 = "  +1.32%  "
$ws.Range("D16").Value2 = "'0.9995"
$ws.Range("D17").Value2 = "'26.528.56"
$ws.Range("E17").Value2 = "  +1.00%  "
$ws.Range("D18").Value2 = "'0.9995"
$ws.Range("E18").Value2 = "  -0.09%  "
$ws.Range("D19").Value2 = "'0.000006943"
$ws.Range("E19").Value2 = "  +2.00%  "
$ws.Range("D20").Value2 = "'11.54"
$ws.Range("E20").Value2 = "  +0.07%  "
$ws.Range("D21").Value2 = "'1.957.13"
$ws.Range("E21").Value2 = "  +1.01%  "
$ws.Range("D22").Value2 = "'4.523"
$ws.Range("E22").Value2 = "  -0.25%  "
$ws.Range("D23").Value2 = "'8.798"
$ws.Range("D24").Value2 = "'5.251"
$ws.Range("E24").Value2 = "  -0.42%  "
$ws.Range("D25").Value2 = "'136.94"
$ws.Range("E25").Value2 = "  +1.51%  "
$ws.Range("E26").Value2 = "  +1.16%  "
$ws.Range("D27").Value2 = "'1.778"
$ws.Range("E27").Value2 = "  +1.08%  "
$ws.Range("E28").Value2 = "  -0.22%  "
$ws.Range("E29").Value2 = "  -0.18%  "
$ws.Range("D30").Value2 = "'3.962"
$ws.Range("E30").Value2 = "  -0.08%  "
$ws.Range("D31").Value2 = "'0.08020"
$ws.Range("E31").Value2 = "  +3.43%  "
$ws.Range("D32").Value2 = "'3.695"
$ws.Range("E32").Value2 = "  +0.53%  "
$ws.Range("D33").Value2 = "'0.04519"
$ws.Range("E33").Value2 = "  +1.41%  "
$ws.Range("D34").Value2 = "'2.612"
$ws.Range("E34").Value2 = "  -0.13%  "
$ws.Range("D35").Value2 = "'0.9985"
$ws.Range("E35").Value2 = "  +2.50%  "
$ws.Range("D36").Value2 = "'0.6262"
$ws.Range("E36").Value2 = "  +1.40%  "
$ws.Range("B37").Value2 = "TrustWalletToken"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value2 = "'0.9118"
$ws.Range("E37").Value2 = "  -1.34%  "
$ws.Range("B38").Value2 = "RenderToken"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value2 = "'2.077"
$ws.Range("E38").Value2 = "  +8.12%  "
$ws.Range("D39").Value2 = "'2.369"
$ws.Range("E39").Value2 = "  -2.53%  "
$ws.Range("D40").Value2 = "'1.001"
$ws.Range("E40").Value2 = "  +0.06%  "
$ws.Range("D41").Value2 = "'103.19"
$ws.Range("E41").Value2 = "  -9.47%  "
$ws.Range("D42").Value2 = "'0.01504"
$ws.Range("E42").Value2 = "  +1.56%  "
$ws.Range("D43").Value2 = "'5.636"
$ws.Range("E43").Value2 = "  +0.41%  "
$ws.Range("D44").Value2 = "'0.3867"
$ws.Range("E44").Value2 = "  +1.21%  "
$ws.Range("D45").Value2 = "'6.978"
$ws.Range("E45").Value2 = "  +11.03%  "
$ws.Range("E46").Value2 = "  +0.20%  "
$ws.Range("D47").Value2 = "'0.05361"
$ws.Range("E47").Value2 = "  +1.73%  "
$ws.Range("D48").Value2 = "'7.828"
$ws.Range("E48").Value2 = "  +0.73%  "
$ws.Range("D49").Value2 = "'30.47"
$ws.Range("E49").Value2 = "  +0.49%  "
$ws.Range("E50").Value2 = "  +3.30%  "
$ws.Range("D51").Value2 = "'51.30"
$ws.Range("E51").Value2 = "  +1.37%  "
